$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# First, clear the stale rows from the old layout. The totals row
# moves from 16 -> 19, and the second "PANE" table moves from rows
# 17-25 -> rows 27-35 (with point values added), so the old row
# positions must be emptied before the new ones are written.
# ------------------------------------------------------------------
$ws.Range("B16:E16").ClearContents()
$ws.Range("A17:G25").ClearContents()

# --- Updates to existing rows (person-on-task / completion reassignments) ---

# Row 5: Admin/Account Status -> person reassigned Vinny -> Alex
$ws.Range("F5").Value = "Alex"

# Row 6: Login -> person Jon -> Alex, completion null -> Done
$ws.Range("F6").Value = "Alex"
$ws.Range("G6").Value = "Done"

# Row 8: Cart -> completion null -> In Progress
$ws.Range("G8").Value = "In Progress"

# Row 9: Search(Interface) -> person assigned Vinny, completion null -> In Progress
$ws.Range("F9").Value = "Vinny"
$ws.Range("G9").Value = "In Progress"

# Row 11: Order/ Date and Time -> completion In Progress -> Done
$ws.Range("G11").Value = "Done"

# Row 12: Order History -> completion In Progress -> Done
$ws.Range("G12").Value = "Done"

# Row 14: Account Info -> completion null -> In Progress
$ws.Range("G14").Value = "In Progress"

# --- New row 15: Database / financials task ---
$ws.Range("A15").Value = "Database / financials "
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 17
$ws.Range("F15").Value = "Vinny"
$ws.Range("G15").Value = "null"

# --- Totals row moves from row 16 to row 19, now including E15 ---
$ws.Range("B19").Value = "1 Point"
$ws.Range("C19").Value = "2 Points"
$ws.Range("D19").Value = "5 Points"
$ws.Range("E19").Formula = "=E2+E3+E4+E5+E6+E7+E8+E9+E10+E11+E12+E13+E14+E15"

# --- Second "PANE" table, moved from rows 17-25 to rows 27-35, with point values added ---
$ws.Range("A28").Value = "Login"
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = "Alex"
$ws.Range("G28").Value = "Done"

$ws.Range("A29").Value = "Register"
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = "Alex"
$ws.Range("G29").Value = "Done"

$ws.Range("A30").Value = "Reset Password/Forgot Email"
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = "Alex"
$ws.Range("G30").Value = "done"

$ws.Range("A27").Value = "PANE (5pnts EACH PANE)"

$ws.Range("A31").Value = "home"
$ws.Range("E31").Value = 10

$ws.Range("A32").Value = "cart/checkout"
$ws.Range("E32").Value = 10

$ws.Range("A33").Value = "display "
$ws.Range("E33").Value = 10

$ws.Range("A34").Value = "information/invoice"
$ws.Range("E34").Value = 10

$ws.Range("A35").Value = "financials"
$ws.Range("E35").Value = 10

# --- Selection: last active cell in the new layout ---
$ws.Range("F23").Select()
